$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 values
$ws.Range("B2").Value = 89.128157874044362
$ws.Range("C2").Value = 61.478991029857156
$ws.Range("D2").Value = 53.494888727437107
$ws.Range("E2").Value = 57.138090646328997

# Update row 3 values
$ws.Range("B3").Value = 74.317296857603409
$ws.Range("C3").Value = 47.41038808743189
$ws.Range("D3").Value = 46.774792690847164
$ws.Range("E3").Value = 57.619929157500792

# Update selection to match the new selected range
$ws.Range("B1:E3").Select()
